$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Swap shared-string order for rows 198 and 199 ---
# Row 198 becomes "San Cristobal y Nieves", Row 199 becomes "Burundi"
$ws.Range("A198").Value = "San Cristobal y Nieves"
$ws.Range("A199").Value = "Burundi"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1161109
$ws.Range("C4").Value = 335
$ws.Range("E4").Value = 919936

# --- Row 9: Alemania ---
$ws.Range("B9").Value = 165016
$ws.Range("C9").Value = 49
$ws.Range("E9").Value = 27604

# --- Row 19: India ---
$ws.Range("B19").Value = 40263
$ws.Range("C19").Value = 564
$ws.Range("D19").Value = 10886
$ws.Range("E19").Value = 28054

# --- Row 22: Arabia Saudita ---
$ws.Range("B22").Value = 27011
$ws.Range("C22").Value = 1552
$ws.Range("D22").Value = 4134
$ws.Range("E22").Value = 22693
$ws.Range("G22").Value = 8
$ws.Range("H22").Value = 184

# --- Row 82: Republica de Macedonia ---
$ws.Range("B82").Value = 1511
$ws.Range("C82").Value = 5
$ws.Range("D82").Value = 945
$ws.Range("E82").Value = 482
$ws.Range("G82").Value = 2
$ws.Range("H82").Value = 84

# --- Row 144: Brunei ---
$ws.Range("D144").Value = 128
$ws.Range("E144").Value = 9

# --- Row 198: San Cristobal y Nieves (new data after string swap) ---
$ws.Range("D198").Value = 8
$ws.Range("H198").Value = 0

# --- Row 199: Burundi (new data after string swap) ---
$ws.Range("D199").Value = 7
$ws.Range("H199").Value = 1
